$wb = $excel.ActiveWorkbook
$metadata = $wb.Worksheets.Item("Metadata")
$concepts = $wb.Worksheets.Item("Concepts")

# --- Metadata sheet (Property / Value table) ---
$metadata.Range("B3").Value = "0.1.2"
$metadata.Range("B5").Value = "CodeSystem - Blood Group (Rh) - NMDP"
$metadata.Range("B8").Value = "2025-04-15T15:35:56-05:00"

# "Count" (B22) holds a text value "2" (was "5"); Value = "2" would be
# auto-coerced to a number by Excel, so force text via NumberFormat then
# restore the original (General) formatting via a format-only paste from a
# sibling cell so the stored style index is unaffected.
$metadata.Range("B22").NumberFormat = "@"
$metadata.Range("B22").Value = "2"
$metadata.Range("B21").Copy() | Out-Null
$metadata.Range("B22").PasteSpecial(-4122) | Out-Null

# --- Concepts sheet (Level / Code / Display / Definition table) ---
# Row 2: I/Indeterminant -> Rh+/Positive
$concepts.Range("B2").Value = "Rh+"
$concepts.Range("C2").Value = "Positive"

# Row 3: P/Positive -> Rh-/Negative
$concepts.Range("B3").Value = "Rh-"
$concepts.Range("C3").Value = "Negative"

# Rows 4-6 (N/Negative, D/Indeterminant, nil/no data) are removed entirely.
$concepts.Range("A4:D6").EntireRow.Delete()
